# Updated cryptos list on Sun Oct 22 08:15:07 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.900.59"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.636.31"
$ws.Range("E3").Value = "  +1.90%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.15%  "

# Row 5 - BNB (force text so "214.94" is not coerced into a number)
$ws.Range("D5").Value = "'214.94"

# Row 6 - XRP
$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.11%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'28.93"
$ws.Range("E8").Value = "  +3.58%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.259"
$ws.Range("E9").Value = "  +2.26%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0610"
$ws.Range("E10").Value = "  +0.93%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0914"
$ws.Range("E11").Value = "  +0.50%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.869.50"
$ws.Range("E12").Value = "  +1.84%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.632.21"
$ws.Range("E13").Value = "  +1.65%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +3.30%  "

# Row 15 - Chainlink
$ws.Range("E15").Value = "  +19.11%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +3.18%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "29.911.34"
$ws.Range("E17").Value = "  +0.76%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "'64.44"
$ws.Range("E18").Value = "  +0.59%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "'243.05"
$ws.Range("E19").Value = "  +0.74%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0703"
$ws.Range("E20").Value = "  +0.62%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.14%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "'9.96"
$ws.Range("E22").Value = "  +5.97%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +3.06%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.89%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'158.12"
$ws.Range("E25").Value = "  +1.75%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'15.62"
$ws.Range("E26").Value = "  +0.95%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +1.37%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +2.47%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.15%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +1.44%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +4.98%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.39"
$ws.Range("E32").Value = "  +4.47%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "'3.18"
$ws.Range("E33").Value = "  -0.19%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.434.47"
$ws.Range("E34").Value = "  +0.47%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "'1.65"
$ws.Range("E35").Value = "  +5.29%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +0.80%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "'2.81"
$ws.Range("E37").Value = "  -4.32%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  +0.08%  "

# Row 39 - Aave
$ws.Range("D39").Value = "'76.84"
$ws.Range("E39").Value = "  +15.68%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "'0.0170"
$ws.Range("E40").Value = "  +0.52%  "

# Row 41 - ImmutableX
$ws.Range("E41").Value = "  +1.12%  "

# Row 42 - RenderToken
$ws.Range("E42").Value = "  +2.21%  "

# Row 43 - ARBITRUM
$ws.Range("E43").Value = "  +1.75%  "

# Row 44 - Kaspa
$ws.Range("D44").Value = "'0.0494"
$ws.Range("E44").Value = "  -1.36%  "

# Rows 45 & 46 swap places: WEMIXToken <-> BitcoinSV, with updated values
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").Value = "'53.78"
$ws.Range("E45").Value = "  -5.42%  "

$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "'1.02"
$ws.Range("E46").Value = "  +2.77%  "

# Row 47 - PaxDollar
$ws.Range("E47").Value = "  +0.14%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.776.00"
$ws.Range("E48").Value = "  +1.90%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "'5.33"
$ws.Range("E49").Value = "  -0.14%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'89.29"
$ws.Range("E50").Value = "  +2.94%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  +5.76%  "
